$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.19"
$ws.Range("D3").Value = "'23.13"
$ws.Range("D4").Value = "'5.392"
$ws.Range("D5").Value = "'0.05928"
$ws.Range("D7").Value = "'0.8067"
$ws.Range("D8").Value = "'0.9100"
$ws.Range("D9").Value = "'0.1416"
$ws.Range("D10").Value = "'0.07438"
$ws.Range("D11").Value = "'0.03326"
$ws.Range("D12").Value = "'0.03044"
$ws.Range("D13").Value = "'0.09325"
$ws.Range("D14").Value = "'3.938"
$ws.Range("D15").Value = "'0.001575"
$ws.Range("D16").Value = "'0.04804"
$ws.Range("D17").Value = "'0.0005942"
$ws.Range("D18").Value = "'0.006107"
$ws.Range("D20").Value = "'0.004420"
$ws.Range("D21").Value = "'0.0009832"
$ws.Range("D22").Value = "'0.00007805"
$ws.Range("D23").Value = "'3.618"
$ws.Range("D24").Value = "'6.445"
$ws.Range("D40").Value = "'0.03875"
$ws.Range("D41").Value = "'0.006212"
$ws.Range("D42").Value = "'0.1065"
$ws.Range("D43").Value = "'0.002802"
$ws.Range("D44").Value = "'0.007237"
$ws.Range("D45").Value = "'0.00005167"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'0.0005802"
$ws.Range("D48").Value = "'0.9804"
